# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, [int]$row, [object[]]$values)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# Sheet "Home win": 7 data rows -> 3 data rows (header + rows 2-4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")

Set-Row $ws 2 @("22-01-2025 22:00", "BRAZIL", "GAÚCHO - 1", "São Luiz - Avenida", 73.3, 1.83)
Set-Row $ws 3 @("22-01-2025 14:50", "SAUDI-ARABIA", "DIVISION 1", "Al Safa - Al Jubail", 70, 2.38)
Set-Row $ws 4 @("23-01-2025 17:45", "WORLD", "UEFA EUROPA LEAGUE", "AZ Alkmaar - AS Roma", 70, 3)

# Remove the now-stale rows 5-7 (whole-row range delete shifts rows up and
# shrinks the used range / <dimension> automatically).
$ws.Range("A5:F7").Delete()

# ---------------------------------------------------------------------------
# Sheet "Draw": 2 data rows -> 4 data rows (header + rows 2-5)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")

Set-Row $ws 2 @("22-01-2025 18:45", "BRAZIL", "CARIOCA - 1", "CFRJ / Maricá - Nova Iguaçu", 60, 3.1)
Set-Row $ws 3 @("22-01-2025 23:00", "BRAZIL", "GOIANO - 1", "Goiatuba EC - Anápolis", 66.7, 2.95)
Set-Row $ws 4 @("23-01-2025 23:00", "BRAZIL", "SERGIPANO", "Barra SE - Lagarto", 60, 5.25)
Set-Row $ws 5 @("23-01-2025 00:00", "NICARAGUA", "PRIMERA DIVISION", "Real Estelí - Walter Ferretti", 73.3, 3.5)

# ---------------------------------------------------------------------------
# Sheet "Btts": 7 data rows -> 4 data rows (header + rows 2-5)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")

Set-Row $ws 2 @("22-01-2025 20:00", "WORLD", "UEFA CHAMPIONS LEAGUE", "Arsenal - Dinamo Zagreb", 76, 2.38)
Set-Row $ws 3 @("22-01-2025 22:00", "COSTA-RICA", "PRIMERA DIVISIÓN", "Municipal Liberia - Santa Ana", 75.59999999999999, 1.77)
Set-Row $ws 4 @("22-01-2025 14:00", "EGYPT", "PREMIER LEAGUE", "Pharco - Al Ahly", 76.7, 2.1)
Set-Row $ws 5 @("14-01-2025 19:45", "FRANCE", "COUPE DE FRANCE", "Haguenau - Dunkerque", 80, 1.93)

$ws.Range("A6:F8").Delete()

# ---------------------------------------------------------------------------
# Sheet "Over_Under": 8 data rows -> 4 data rows (header + rows 2-5)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")

Set-Row $ws 2 @("22-01-2025 17:45", "WORLD", "UEFA CHAMPIONS LEAGUE", "Shakhtar Donetsk - Stade Brestois 29", 73.3, 1.91, 66.7, 3.4)
Set-Row $ws 3 @("23-01-2025 17:45", "WORLD", "UEFA EUROPA LEAGUE", "Bodo/Glimt - Maccabi Tel Aviv", 93.3, 1.62, 66.7, 2.5)
Set-Row $ws 4 @("23-01-2025 17:45", "WORLD", "UEFA EUROPA LEAGUE", "Malmo FF - Twente", 73.3, 1.73, 60, 2.75)
Set-Row $ws 5 @("23-01-2025 20:00", "WORLD", "UEFA EUROPA LEAGUE", "Manchester United - Rangers", 80, 1.62, 60, 2.5)

$ws.Range("A6:H9").Delete()
